$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 251
$ws1.Range("F3").Value = 2484
$ws1.Range("F5").Value = 904
$ws1.Range("F7").Value = 1305
$ws1.Range("F8").Value = 1660
$ws1.Range("F10").Value = 59
$ws1.Range("F11").Value = 2322
$ws1.Range("F12").Value = 475
$ws1.Range("F13").Value = 148
$ws1.Range("F18").Value = 8455
$ws1.Range("F20").Value = 6517
$ws1.Range("F21").Value = 10525
$ws1.Range("F24").Value = 200
$ws1.Range("F25").Value = 289
$ws1.Range("F26").Value = 518
$ws1.Range("F27").Value = 187
$ws1.Range("F28").Value = 165
$ws1.Range("F29").Value = 72
$ws1.Range("F30").Value = 16
$ws1.Range("F31").Value = 10
$ws1.Range("F32").Value = 4431
$ws1.Range("F33").Value = 321
$ws1.Range("F34").Value = 420

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F16").Value = 93

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 612

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 612
$ws4.Range("F4").Value = 251
$ws4.Range("F6").Value = 2484
$ws4.Range("F8").Value = 904
$ws4.Range("F10").Value = 1305
$ws4.Range("F12").Value = 1660
$ws4.Range("F15").Value = 2322
$ws4.Range("F17").Value = 475
$ws4.Range("F18").Value = 148
$ws4.Range("F24").Value = 8455
$ws4.Range("F26").Value = 6517
$ws4.Range("F27").Value = 10525
$ws4.Range("F31").Value = 200
$ws4.Range("F32").Value = 289
$ws4.Range("F34").Value = 518
$ws4.Range("F38").Value = 187
$ws4.Range("F39").Value = 165
$ws4.Range("F40").Value = 4431
$ws4.Range("F43").Value = 93
$ws4.Range("F47").Value = 420
